$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("PIR")

$ws.Range("A59:A72").NumberFormat = "@"

$ws.Range("A59").Value = "2026-01-28"
$ws.Range("B59").Value = "17:06:50"
$ws.Range("C59").Value = "17:00"
$ws.Range("D59").Value = "Bathroom"
$ws.Range("E59").Value = "No Motion"
$ws.Range("F59").Value = "Inactive"
$ws.Range("A60").Value = "2026-01-28"
$ws.Range("B60").Value = "17:06:51"
$ws.Range("C60").Value = "17:00"
$ws.Range("D60").Value = "Bathroom"
$ws.Range("E60").Value = "No Motion"
$ws.Range("F60").Value = "Inactive"
$ws.Range("A61").Value = "2026-01-28"
$ws.Range("B61").Value = "17:06:54"
$ws.Range("C61").Value = "17:00"
$ws.Range("D61").Value = "Bathroom"
$ws.Range("E61").Value = "No Motion"
$ws.Range("F61").Value = "Inactive"
$ws.Range("A62").Value = "2026-01-28"
$ws.Range("B62").Value = "17:06:58"
$ws.Range("C62").Value = "17:00"
$ws.Range("D62").Value = "Bathroom"
$ws.Range("E62").Value = "No Motion"
$ws.Range("F62").Value = "Inactive"
$ws.Range("A63").Value = "2026-01-28"
$ws.Range("B63").Value = "17:07:04"
$ws.Range("C63").Value = "17:00"
$ws.Range("D63").Value = "Bathroom"
$ws.Range("E63").Value = "No Motion"
$ws.Range("F63").Value = "Inactive"
$ws.Range("A64").Value = "2026-01-28"
$ws.Range("B64").Value = "17:07:09"
$ws.Range("C64").Value = "17:00"
$ws.Range("D64").Value = "Bathroom"
$ws.Range("E64").Value = "No Motion"
$ws.Range("F64").Value = "Inactive"
$ws.Range("A65").Value = "2026-01-28"
$ws.Range("B65").Value = "17:07:14"
$ws.Range("C65").Value = "17:00"
$ws.Range("D65").Value = "Bathroom"
$ws.Range("E65").Value = "No Motion"
$ws.Range("F65").Value = "Inactive"
$ws.Range("A66").Value = "2026-01-28"
$ws.Range("B66").Value = "17:07:19"
$ws.Range("C66").Value = "17:00"
$ws.Range("D66").Value = "Bathroom"
$ws.Range("E66").Value = "No Motion"
$ws.Range("F66").Value = "Inactive"
$ws.Range("A67").Value = "2026-01-28"
$ws.Range("B67").Value = "17:07:24"
$ws.Range("C67").Value = "17:00"
$ws.Range("D67").Value = "Bathroom"
$ws.Range("E67").Value = "No Motion"
$ws.Range("F67").Value = "Inactive"
$ws.Range("A68").Value = "2026-01-28"
$ws.Range("B68").Value = "17:07:29"
$ws.Range("C68").Value = "17:00"
$ws.Range("D68").Value = "Bathroom"
$ws.Range("E68").Value = "No Motion"
$ws.Range("F68").Value = "Inactive"
$ws.Range("A69").Value = "2026-01-28"
$ws.Range("B69").Value = "17:07:34"
$ws.Range("C69").Value = "17:00"
$ws.Range("D69").Value = "Bathroom"
$ws.Range("E69").Value = "No Motion"
$ws.Range("F69").Value = "Inactive"
$ws.Range("A70").Value = "2026-01-28"
$ws.Range("B70").Value = "17:07:39"
$ws.Range("C70").Value = "17:00"
$ws.Range("D70").Value = "Bathroom"
$ws.Range("E70").Value = "No Motion"
$ws.Range("F70").Value = "Inactive"
$ws.Range("A71").Value = "2026-01-28"
$ws.Range("B71").Value = "17:07:44"
$ws.Range("C71").Value = "17:00"
$ws.Range("D71").Value = "Bathroom"
$ws.Range("E71").Value = "No Motion"
$ws.Range("F71").Value = "Inactive"
$ws.Range("A72").Value = "2026-01-28"
$ws.Range("B72").Value = "17:07:49"
$ws.Range("C72").Value = "17:00"
$ws.Range("D72").Value = "Bathroom"
$ws.Range("E72").Value = "No Motion"
$ws.Range("F72").Value = "Inactive"

$ws = $wb.Worksheets.Item("Humidity")

$ws.Range("A59:A72").NumberFormat = "@"
$ws.Range("E59:E72").NumberFormat = "@"

$ws.Range("A59").Value = "2026-01-28"
$ws.Range("B59").Value = "17:06:50"
$ws.Range("C59").Value = "17:00"
$ws.Range("D59").Value = "Bathroom"
$ws.Range("E59").Value = "86.7%"
$ws.Range("F59").Value = "Active"
$ws.Range("A60").Value = "2026-01-28"
$ws.Range("B60").Value = "17:06:50"
$ws.Range("C60").Value = "17:00"
$ws.Range("D60").Value = "Bathroom"
$ws.Range("E60").Value = "87.6%"
$ws.Range("F60").Value = "Active"
$ws.Range("A61").Value = "2026-01-28"
$ws.Range("B61").Value = "17:06:56"
$ws.Range("C61").Value = "17:00"
$ws.Range("D61").Value = "Bathroom"
$ws.Range("E61").Value = "86.7%"
$ws.Range("F61").Value = "Active"
$ws.Range("A62").Value = "2026-01-28"
$ws.Range("B62").Value = "17:07:00"
$ws.Range("C62").Value = "17:00"
$ws.Range("D62").Value = "Bathroom"
$ws.Range("E62").Value = "87.6%"
$ws.Range("F62").Value = "Active"
$ws.Range("A63").Value = "2026-01-28"
$ws.Range("B63").Value = "17:07:04"
$ws.Range("C63").Value = "17:00"
$ws.Range("D63").Value = "Bathroom"
$ws.Range("E63").Value = "86.6%"
$ws.Range("F63").Value = "Active"
$ws.Range("A64").Value = "2026-01-28"
$ws.Range("B64").Value = "17:07:16"
$ws.Range("C64").Value = "17:00"
$ws.Range("D64").Value = "Bathroom"
$ws.Range("E64").Value = "87.5%"
$ws.Range("F64").Value = "Active"
$ws.Range("A65").Value = "2026-01-28"
$ws.Range("B65").Value = "17:07:20"
$ws.Range("C65").Value = "17:00"
$ws.Range("D65").Value = "Bathroom"
$ws.Range("E65").Value = "87.6%"
$ws.Range("F65").Value = "Active"
$ws.Range("A66").Value = "2026-01-28"
$ws.Range("B66").Value = "17:07:24"
$ws.Range("C66").Value = "17:00"
$ws.Range("D66").Value = "Bathroom"
$ws.Range("E66").Value = "86.7%"
$ws.Range("F66").Value = "Active"
$ws.Range("A67").Value = "2026-01-28"
$ws.Range("B67").Value = "17:07:28"
$ws.Range("C67").Value = "17:00"
$ws.Range("D67").Value = "Bathroom"
$ws.Range("E67").Value = "87.5%"
$ws.Range("F67").Value = "Active"
$ws.Range("A68").Value = "2026-01-28"
$ws.Range("B68").Value = "17:07:32"
$ws.Range("C68").Value = "17:00"
$ws.Range("D68").Value = "Bathroom"
$ws.Range("E68").Value = "87.5%"
$ws.Range("F68").Value = "Active"
$ws.Range("A69").Value = "2026-01-28"
$ws.Range("B69").Value = "17:07:36"
$ws.Range("C69").Value = "17:00"
$ws.Range("D69").Value = "Bathroom"
$ws.Range("E69").Value = "86.6%"
$ws.Range("F69").Value = "Active"
$ws.Range("A70").Value = "2026-01-28"
$ws.Range("B70").Value = "17:07:40"
$ws.Range("C70").Value = "17:00"
$ws.Range("D70").Value = "Bathroom"
$ws.Range("E70").Value = "87.6%"
$ws.Range("F70").Value = "Active"
$ws.Range("A71").Value = "2026-01-28"
$ws.Range("B71").Value = "17:07:44"
$ws.Range("C71").Value = "17:00"
$ws.Range("D71").Value = "Bathroom"
$ws.Range("E71").Value = "86.6%"
$ws.Range("F71").Value = "Active"
$ws.Range("A72").Value = "2026-01-28"
$ws.Range("B72").Value = "17:07:48"
$ws.Range("C72").Value = "17:00"
$ws.Range("D72").Value = "Bathroom"
$ws.Range("E72").Value = "86.1%"
$ws.Range("F72").Value = "Active"

$ws = $wb.Worksheets.Item("Temperature")

$ws.Range("A59:A72").NumberFormat = "@"

$ws.Range("A59").Value = "2026-01-28"
$ws.Range("B59").Value = "17:06:50"
$ws.Range("C59").Value = "17:00"
$ws.Range("D59").Value = "Bathroom"
$ws.Range("E59").Value = "22.9C"
$ws.Range("F59").Value = "Active"
$ws.Range("A60").Value = "2026-01-28"
$ws.Range("B60").Value = "17:06:51"
$ws.Range("C60").Value = "17:00"
$ws.Range("D60").Value = "Bathroom"
$ws.Range("E60").Value = "22.9C"
$ws.Range("F60").Value = "Active"
$ws.Range("A61").Value = "2026-01-28"
$ws.Range("B61").Value = "17:06:56"
$ws.Range("C61").Value = "17:00"
$ws.Range("D61").Value = "Bathroom"
$ws.Range("E61").Value = "22.9C"
$ws.Range("F61").Value = "Active"
$ws.Range("A62").Value = "2026-01-28"
$ws.Range("B62").Value = "17:07:00"
$ws.Range("C62").Value = "17:00"
$ws.Range("D62").Value = "Bathroom"
$ws.Range("E62").Value = "22.8C"
$ws.Range("F62").Value = "Active"
$ws.Range("A63").Value = "2026-01-28"
$ws.Range("B63").Value = "17:07:04"
$ws.Range("C63").Value = "17:00"
$ws.Range("D63").Value = "Bathroom"
$ws.Range("E63").Value = "22.8C"
$ws.Range("F63").Value = "Active"
$ws.Range("A64").Value = "2026-01-28"
$ws.Range("B64").Value = "17:07:16"
$ws.Range("C64").Value = "17:00"
$ws.Range("D64").Value = "Bathroom"
$ws.Range("E64").Value = "22.8C"
$ws.Range("F64").Value = "Active"
$ws.Range("A65").Value = "2026-01-28"
$ws.Range("B65").Value = "17:07:20"
$ws.Range("C65").Value = "17:00"
$ws.Range("D65").Value = "Bathroom"
$ws.Range("E65").Value = "22.9C"
$ws.Range("F65").Value = "Active"
$ws.Range("A66").Value = "2026-01-28"
$ws.Range("B66").Value = "17:07:24"
$ws.Range("C66").Value = "17:00"
$ws.Range("D66").Value = "Bathroom"
$ws.Range("E66").Value = "22.9C"
$ws.Range("F66").Value = "Active"
$ws.Range("A67").Value = "2026-01-28"
$ws.Range("B67").Value = "17:07:28"
$ws.Range("C67").Value = "17:00"
$ws.Range("D67").Value = "Bathroom"
$ws.Range("E67").Value = "22.8C"
$ws.Range("F67").Value = "Active"
$ws.Range("A68").Value = "2026-01-28"
$ws.Range("B68").Value = "17:07:32"
$ws.Range("C68").Value = "17:00"
$ws.Range("D68").Value = "Bathroom"
$ws.Range("E68").Value = "22.8C"
$ws.Range("F68").Value = "Active"
$ws.Range("A69").Value = "2026-01-28"
$ws.Range("B69").Value = "17:07:36"
$ws.Range("C69").Value = "17:00"
$ws.Range("D69").Value = "Bathroom"
$ws.Range("E69").Value = "22.8C"
$ws.Range("F69").Value = "Active"
$ws.Range("A70").Value = "2026-01-28"
$ws.Range("B70").Value = "17:07:40"
$ws.Range("C70").Value = "17:00"
$ws.Range("D70").Value = "Bathroom"
$ws.Range("E70").Value = "22.9C"
$ws.Range("F70").Value = "Active"
$ws.Range("A71").Value = "2026-01-28"
$ws.Range("B71").Value = "17:07:44"
$ws.Range("C71").Value = "17:00"
$ws.Range("D71").Value = "Bathroom"
$ws.Range("E71").Value = "22.8C"
$ws.Range("F71").Value = "Active"
$ws.Range("A72").Value = "2026-01-28"
$ws.Range("B72").Value = "17:07:48"
$ws.Range("C72").Value = "17:00"
$ws.Range("D72").Value = "Bathroom"
$ws.Range("E72").Value = "22.8C"
$ws.Range("F72").Value = "Active"
